# Fix emails in the worksheet and add a hyperlink on G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the e-mail strings that had the number appended at the wrong place.
$ws.Range("G4").Value = "nemreg1es1@mail.com"
$ws.Range("G5").Value = "nemreg1es2@mail.com"
$ws.Range("G6").Value = "nemreg1es3@mail.com"

# These two cells held mistaken e-mail addresses that must be cleared entirely.
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Add the hyperlink for the corrected G5 e-mail address. Hyperlinks.Add
# applies Excel's built-in "Hyperlink" cell style automatically, so copy the
# plain formatting back from a neighbouring cell to keep G5 looking the same
# as before.
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:nemreg1es2@mail.com", "", "", "nemreg1es2@mail.com")
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)

# Row 5 height change.
$ws.Rows.Item(5).RowHeight = 15.65

# Update the active selection to match the authored state.
$ws.Range("G6").Select()
